{"js": "// Update the two-digit \u00f7 one-digit division answer table cells to the\n// newly generated set of problems (commit \"Update master to output\n// generated at c8c62b6\"). Each old cell string is unique in the\n// document, so a direct search+replace per pair is safe and will not\n// disturb the one cell that stays the same (63\u00f79=7, 0).\nconst replacements = [\n  [\"17\u00f79=1, 8\", \"37\u00f74=9, 1\"],\n  [\"96\u00f75=19, 1\", \"37\u00f79=4, 1\"],\n  [\"31\u00f74=7, 3\", \"72\u00f78=9, 0\"],\n  [\"56\u00f79=6, 2\", \"92\u00f75=18, 2\"],\n  [\"26\u00f72=13, 0\", \"27\u00f72=13, 1\"],\n  [\"38\u00f73=12, 2\", \"56\u00f73=18, 2\"],\n  [\"38\u00f79=4, 2\", \"74\u00f76=12, 2\"],\n  [\"24\u00f76=4, 0\", \"10\u00f76=1, 4\"],\n  [\"45\u00f76=7, 3\", \"95\u00f72=47, 1\"],\n  [\"37\u00f77=5, 2\", \"40\u00f74=10, 0\"],\n  [\"19\u00f78=2, 3\", \"35\u00f76=5, 5\"],\n  [\"60\u00f76=10, 0\", \"73\u00f76=12, 1\"],\n  [\"92\u00f76=15, 2\", \"18\u00f78=2, 2\"],\n  [\"33\u00f72=16, 1\", \"79\u00f73=26, 1\"],\n  [\"54\u00f72=27, 0\", \"70\u00f75=14, 0\"],\n  [\"77\u00f78=9, 5\", \"54\u00f79=6, 0\"],\n  [\"62\u00f75=12, 2\", \"50\u00f72=25, 0\"],\n  [\"76\u00f75=15, 1\", \"54\u00f79=6, 0\"],\n  [\"53\u00f75=10, 3\", \"96\u00f79=10, 6\"],\n  [\"96\u00f78=12, 0\", \"93\u00f77=13, 2\"],\n  [\"65\u00f72=32, 1\", \"96\u00f73=32, 0\"],\n  [\"49\u00f76=8, 1\", \"88\u00f76=14, 4\"],\n  [\"20\u00f73=6, 2\", \"46\u00f75=9, 1\"],\n  [\"87\u00f72=43, 1\", \"56\u00f72=28, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit / one-digit division answer table cells to the\n# newly generated set of problems (commit \"Update master to output\n# generated at c8c62b6\"). Each old cell string is unique in the\n# document, so Find/Replace per pair is safe and leaves the one\n# unchanged cell (63\u00f79=7, 0) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old = \"17\u00f79=1, 8\";   new = \"37\u00f74=9, 1\"},\n  @{old = \"96\u00f75=19, 1\";  new = \"37\u00f79=4, 1\"},\n  @{old = \"31\u00f74=7, 3\";   new = \"72\u00f78=9, 0\"},\n  @{old = \"56\u00f79=6, 2\";   new = \"92\u00f75=18, 2\"},\n  @{old = \"26\u00f72=13, 0\";  new = \"27\u00f72=13, 1\"},\n  @{old = \"38\u00f73=12, 2\";  new = \"56\u00f73=18, 2\"},\n  @{old = \"38\u00f79=4, 2\";   new = \"74\u00f76=12, 2\"},\n  @{old = \"24\u00f76=4, 0\";   new = \"10\u00f76=1, 4\"},\n  @{old = \"45\u00f76=7, 3\";   new = \"95\u00f72=47, 1\"},\n  @{old = \"37\u00f77=5, 2\";   new = \"40\u00f74=10, 0\"},\n  @{old = \"19\u00f78=2, 3\";   new = \"35\u00f76=5, 5\"},\n  @{old = \"60\u00f76=10, 0\";  new = \"73\u00f76=12, 1\"},\n  @{old = \"92\u00f76=15, 2\";  new = \"18\u00f78=2, 2\"},\n  @{old = \"33\u00f72=16, 1\";  new = \"79\u00f73=26, 1\"},\n  @{old = \"54\u00f72=27, 0\";  new = \"70\u00f75=14, 0\"},\n  @{old = \"77\u00f78=9, 5\";   new = \"54\u00f79=6, 0\"},\n  @{old = \"62\u00f75=12, 2\";  new = \"50\u00f72=25, 0\"},\n  @{old = \"76\u00f75=15, 1\";  new = \"54\u00f79=6, 0\"},\n  @{old = \"53\u00f75=10, 3\";  new = \"96\u00f79=10, 6\"},\n  @{old = \"96\u00f78=12, 0\";  new = \"93\u00f77=13, 2\"},\n  @{old = \"65\u00f72=32, 1\";  new = \"96\u00f73=32, 0\"},\n  @{old = \"49\u00f76=8, 1\";   new = \"88\u00f76=14, 4\"},\n  @{old = \"20\u00f73=6, 2\";   new = \"46\u00f75=9, 1\"},\n  @{old = \"87\u00f72=43, 1\";  new = \"56\u00f72=28, 0\"}\n)\n\nforeach ($p in $pairs) {\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Text = $p.old\n  $find.Replacement.Text = $p.new\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)\n}\n"}
